# R_Commands.xlsx edit:
# - Insert a new "Restart R session" / ".rs.restartR()" row into the
#   alphabetically-sorted list of R commands (between "Render R markdown"
#   and "Run script"), which lands at row 15 and pushes later rows down.
# - Leave the active selection on B11 (the Install package syntax cell).
# - Nudge the page setup (portrait orientation) as part of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 15 ("Run script") so the new
# "Restart R session" entry keeps the list in alphabetical order.
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = "Restart R session"
$ws.Range("B15").Value = ".rs.restartR()"

# Restore/confirm page orientation (portrait) for the worksheet.
$ws.PageSetup.Orientation = 1

# Leave selection on B11, matching the saved view state.
$ws.Range("B11").Select()
